# Insert a new price-report row for "Apio" (Vega Modelo de Temuco) ahead of
# the existing row 421, shifting all subsequent rows down by one (421-481
# become 422-482). The new row 421 carries a fresh weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 421; Excel shifts rows 421..481 down to 422..482
# and the used range / dimension grows to row 482 automatically.
$ws.Rows.Item(421).Insert()

# Populate the newly inserted row 421 with the new observation.
$ws.Cells.Item(421, 1).Value = 10
$ws.Cells.Item(421, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(421, 3).Value = "La Araucanía"
$ws.Cells.Item(421, 4).Value = 45077
$ws.Cells.Item(421, 5).Value = 9
$ws.Cells.Item(421, 6).Value = 100112017
$ws.Cells.Item(421, 7).Value = "Apio"
$ws.Cells.Item(421, 8).Value = "Americana (o)"
$ws.Cells.Item(421, 9).Value = "Primera"
$ws.Cells.Item(421, 10).Value = 55
$ws.Cells.Item(421, 11).Value = 8000
$ws.Cells.Item(421, 12).Value = 8000
$ws.Cells.Item(421, 13).Value = 8000
$ws.Cells.Item(421, 14).Value = "$/docena de matas"
$ws.Cells.Item(421, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(421, 16).Value = 1333
$ws.Cells.Item(421, 17).Value = 6
$ws.Cells.Item(421, 18).Value = "Hortaliza"
